# Atualizado por script em 27-11-2023 20:21
#
# Applies the betexplorer scrape refresh:
#   - 7 pairs of adjacent rows (sharing the same kickoff datetime) get
#     re-ordered, i.e. columns F:V are swapped between the two rows.
#   - 2 new match rows are appended at the end (rows 130 and 131).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap pairs (same E/date, only F:V content reordered) ---
$swapPairs = @(
    @(4, 5),
    @(6, 7),
    @(16, 17),
    @(18, 19),
    @(56, 58),
    @(59, 60),
    @(86, 87)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("F$r1`:V$r1")
    $range2 = $ws.Range("F$r2`:V$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value = $vals2
    $range2.Value = $vals1
}

# --- Append two new match rows at the bottom ---

# Row 130 (Indice 129): Verona x Lecce
$ws.Range("A130").Value = 129
$ws.Range("B130").Value = "italy"
$ws.Range("C130").Value = "serie-a"
$ws.Range("D130").Value = "2023-2024"
$ws.Range("E130").Value = 45257.77083333334

$arr130 = New-Object 'object[,]' 1,17
$arr130[0,0]  = "Verona"
$arr130[0,1]  = 2
$arr130[0,2]  = "Lecce"
$arr130[0,3]  = 2
$arr130[0,4]  = 2.6
$arr130[0,5]  = "05/11/2023 11:03"
$arr130[0,6]  = 2.88
$arr130[0,7]  = "27/11/2023 18:28"
$arr130[0,8]  = 3.04
$arr130[0,9]  = "05/11/2023 11:03"
$arr130[0,10] = 3.02
$arr130[0,11] = "27/11/2023 18:15"
$arr130[0,12] = 3.07
$arr130[0,13] = "05/11/2023 11:03"
$arr130[0,14] = 2.87
$arr130[0,15] = "27/11/2023 18:28"
$arr130[0,16] = "https://www.betexplorer.com/football/italy/serie-a/verona-lecce/nT4T8lPk/"
$ws.Range("F130:V130").Value = $arr130

# Row 131 (Indice 130): Bologna x Torino
$ws.Range("A131").Value = 130
$ws.Range("B131").Value = "italy"
$ws.Range("C131").Value = "serie-a"
$ws.Range("D131").Value = "2023-2024"
$ws.Range("E131").Value = 45257.86458333334

$arr131 = New-Object 'object[,]' 1,17
$arr131[0,0]  = "Bologna"
$arr131[0,1]  = 2
$arr131[0,2]  = "Torino"
$arr131[0,3]  = 0
$arr131[0,4]  = 2.34
$arr131[0,5]  = "05/11/2023 11:03"
$arr131[0,6]  = 2.38
$arr131[0,7]  = "27/11/2023 20:26"
$arr131[0,8]  = 3.17
$arr131[0,9]  = "05/11/2023 11:03"
$arr131[0,10] = 3.07
$arr131[0,11] = "27/11/2023 20:44"
$arr131[0,12] = 3.36
$arr131[0,13] = "05/11/2023 11:03"
$arr131[0,14] = 3.56
$arr131[0,15] = "27/11/2023 20:44"
$arr131[0,16] = "https://www.betexplorer.com/football/italy/serie-a/bologna-torino/8MH1fkAd/"
$ws.Range("F131:V131").Value = $arr131

# --- Match formatting of the new rows to the existing table rows ---
$ws.Range("A130:A131").Font.Bold = $true
$ws.Range("A130:A131").HorizontalAlignment = -4108
$ws.Range("A130:A131").VerticalAlignment = -4160
$ws.Range("A130:A131").Borders.LineStyle = 1

$ws.Range("E130:E131").NumberFormat = $ws.Range("E129").NumberFormat
